# edit.ps1 - applies "added Setup und Implementierung" changes
# to MIDI_Baton_Dokumentation.docx

$d = $word.ActiveDocument

function Replace-Exact {
    param(
        [string]$OldText,
        [string]$NewText
    )
    $range = $d.Content
    $ok = $range.Find.Execute(
        $OldText, $false, $false, $false, $false, $false,
        $true, 1, $false, $NewText, 2
    )
    if (-not $ok) {
        Write-Host "WARNING: replacement failed for: $OldText"
    }
    return $ok
}

# ---------------------------------------------------------------------------
# 1) "Implementierung" section (MIDI-Baton concept chapter): describe the
#    actual built prototype instead of the placeholder text.
# ---------------------------------------------------------------------------
Replace-Exact `
    "Der Prototyp besteht aus einem Holzgerüst, das mit Kreppband bezogen wurde. LED-Streifen lassen die Kanten in wechselnden Regenbogenfarben erstrahlen. ..." `
    "Der Prototyp besteht aus einem Plastikrohr, an dessen Ende ein Beschleunigungssensor unter einem Tischtennisball versteckt ist, im Rohr ist außerdem ein LED-Streifen verbaut. Das Rohr ist am unteren Ende mit Filz umwickelt und darüber ein Drucksensor. Unterhalb des Griffs aus Filz führen die Kabel in eine Box die via Klettverschluss an einem Armband am Handgelenk befestigt werden kann. Darin ist ein Arduino verbaut und die Box verlässt das USB-Kabel, dass zum Anschluss des Batons benötigt wird."

# ---------------------------------------------------------------------------
# 2) "Setup" section body text.
# ---------------------------------------------------------------------------
Replace-Exact `
    "Beschreibung, wie man den Prototypen demonstrationsfähig macht." `
    "Zunächst sollte man sich die Schlaufe mit der Box ums Handgelenk legen, danach schließt man das Gerät an einen MIDI-Host an und sobald es erkannt wird kann gespielt werden."

# ---------------------------------------------------------------------------
# 3) "Bedienungsanleitung" heading: merge the two runs ("Bedienungsanleit" +
#    "ung") into a single run with identical text - use the paragraph
#    object directly, since the same word also appears as a header caption
#    earlier in the document.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Bedienungsanleit" + "ung" -and $p.Style.NameLocal -eq "Heading 3") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Text = "Bedienungsanleitung"
        break
    }
}

# ---------------------------------------------------------------------------
# 4) "Beschreibung, wie man mit dem Prototypen interagieren kann" body text.
# ---------------------------------------------------------------------------
Replace-Exact `
    "Beschreibung, wie man mit dem Prototypen interagieren kann" `
    "Der Prototyp kann via Neigung nach oben oder unten Noten spielen, neigt man ihn höher werden auch die Noten höher und umgekehrt. Um die Note auch abzuspielen muss der Drucksensor über dem Griff betätigt werden, dessen gemessener Druck bestimmt auch die Lautstärke der angeschlagenen Note, wenn der MIDI-Host dies unterstützt."

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark from the "Alexander Eder" paragraph to the
#    page-break paragraph right before the "Dokumentation" section (this
#    simply reflects where the cursor was when the document was last saved).
# ---------------------------------------------------------------------------
$bookmarks = $d.Bookmarks
if ($bookmarks.Exists("_GoBack")) {
    $bookmarks.Item("_GoBack").Delete()
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "*interagieren*" -or $txt -like "*MIDI-Host dies*") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        $target = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)
        $bookmarks.Add("_GoBack", $target) | Out-Null
        break
    }
}
